# Deploy the implementation guide.
#
# The "Metadata" sheet lists the ValueSet's FHIR metadata as Property/Value
# rows. This refreshes the published Status (was "active", now still a
# "draft") and the Date the guide was generated/published.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 6: Status -> draft (was "active")
$ws.Range("B6").Value = "draft"

# Row 8: Date -> new publish timestamp (was 2023-05-12T12:33:13+00:00)
$ws.Range("B8").Value = "2023-08-01T16:12:28+00:00"
